$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:B1): bump font size 11 -> 12, and give row 1 an explicit height ---
$ws.Range("A1:B1").Font.Size = 12
$ws.Rows.Item(1).RowHeight = 15.75

# --- Spacer row (A2:B2): fill with the "Text 1" theme color (black) ---
$ws.Range("A2:B2").Interior.ThemeColor = 1

# --- New row 14: add the new "go for dinner" entry at 19:30 ---
$ws.Range("A14").Value = 0.8125
$ws.Range("A14").NumberFormat = "h:mm AM/PM"
$ws.Range("A14").HorizontalAlignment = -4108
$ws.Range("B14").Value = "Go for the dinner by shutting down pc"

# --- Update the selection / view state ---
$ws.Range("D15").Select()
